# Atualização de bases das ligas, do dia: 23-02-2024 às 08:18
#
# The underlying data feed re-sorted/re-paired a handful of match rows:
#   - Rows 262/263 (ids 6924569 / 6924568) swap their full data (the
#     match that was on row 262 moves to row 263 and vice-versa).
#   - Rows 339/340/341 gain previously-missing FTHG/FTAG/FTR values and
#     have their closing-odds / P&L columns refreshed.
#   - Row 342 (id 7641677) is replaced by what used to be row 343
#     (id 7641680, with slightly updated oddD/oddA), and the old row 343
#     is removed outright (the sheet shrinks from 343 to 342 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 262 (becomes the former row-263 match: Atletico Morelia vs Atlante) ---
$ws.Range("B262").Value = 6924568
$ws.Range("F262").Value = "Atletico Morelia"
$ws.Range("G262").Value = "Atlante"
$ws.Range("H262").Value = 0
$ws.Range("I262").Value = 1
$ws.Range("J262").Value = "A"
$ws.Range("K262").Value = 2.4
$ws.Range("L262").Value = 3
$ws.Range("M262").Value = 2.875
$ws.Range("N262").Value = 2.7
$ws.Range("O262").Value = 3.1
$ws.Range("P262").Value = 2.8
$ws.Range("Q262").Value = 0
$ws.Range("R262").Value = 1.85
$ws.Range("S262").Value = 1.95
$ws.Range("T262").Value = 2.25
$ws.Range("U262").Value = 1.975
$ws.Range("V262").Value = 1.725
$ws.Range("W262").Value = -1
$ws.Range("X262").Value = -1
$ws.Range("Y262").Value = 1.8
$ws.Range("Z262").Value = -1
$ws.Range("AA262").Value = 0.95
$ws.Range("AB262").Value = -1
$ws.Range("AC262").Value = 0.7250000000000001

# --- Row 263 (becomes the former row-262 match: Venados FC vs Dorados) ---
$ws.Range("B263").Value = 6924569
$ws.Range("F263").Value = "Venados FC"
$ws.Range("G263").Value = "Dorados"
$ws.Range("H263").Value = 4
$ws.Range("I263").Value = 1
$ws.Range("J263").Value = "H"
$ws.Range("K263").Value = 1.615
$ws.Range("L263").Value = 4
$ws.Range("M263").Value = 4.5
$ws.Range("N263").Value = 1.5
$ws.Range("O263").Value = 4.75
$ws.Range("P263").Value = 5.75
$ws.Range("Q263").Value = -1.25
$ws.Range("R263").Value = 1.925
$ws.Range("S263").Value = 1.875
$ws.Range("T263").Value = 3
$ws.Range("U263").Value = 1.75
$ws.Range("V263").Value = 1.95
$ws.Range("W263").Value = 0.5
$ws.Range("X263").Value = -1
$ws.Range("Y263").Value = -1
$ws.Range("Z263").Value = 0.925
$ws.Range("AA263").Value = -1
$ws.Range("AB263").Value = 0.75
$ws.Range("AC263").Value = -1

# --- Row 339: add the final-score columns + refresh closing odds / PL ---
$ws.Range("H339").Value = 1
$ws.Range("I339").Value = 0
$ws.Range("J339").Value = "H"
$ws.Range("R339").Value = 2.025
$ws.Range("S339").Value = 1.775
$ws.Range("T339").Value = 2.25
$ws.Range("U339").Value = 1.75
$ws.Range("V339").Value = 1.95
$ws.Range("W339").Value = 1.9
$ws.Range("X339").Value = -1
$ws.Range("Y339").Value = -1
$ws.Range("Z339").Value = 1.025
$ws.Range("AA339").Value = -1
$ws.Range("AB339").Value = -1
$ws.Range("AC339").Value = 0.95

# --- Row 340: add the final-score columns + refresh closing odds / PL ---
$ws.Range("H340").Value = 4
$ws.Range("I340").Value = 0
$ws.Range("J340").Value = "H"
$ws.Range("N340").Value = 1.833
$ws.Range("O340").Value = 3.8
$ws.Range("P340").Value = 4.2
$ws.Range("R340").Value = 1.8
$ws.Range("S340").Value = 2
$ws.Range("T340").Value = 2.75
$ws.Range("U340").Value = 2
$ws.Range("V340").Value = 1.8
$ws.Range("W340").Value = 0.833
$ws.Range("X340").Value = -1
$ws.Range("Y340").Value = -1
$ws.Range("Z340").Value = 0.8
$ws.Range("AA340").Value = -1
$ws.Range("AB340").Value = 1
$ws.Range("AC340").Value = -1

# --- Row 341: add the final-score columns + refresh closing odds / PL ---
$ws.Range("H341").Value = 1
$ws.Range("I341").Value = 0
$ws.Range("J341").Value = "H"
$ws.Range("R341").Value = 1.825
$ws.Range("S341").Value = 1.975
$ws.Range("T341").Value = 2.25
$ws.Range("U341").Value = 1.75
$ws.Range("V341").Value = 1.95
$ws.Range("W341").Value = 0.8500000000000001
$ws.Range("X341").Value = -1
$ws.Range("Y341").Value = -1
$ws.Range("Z341").Value = 0.825
$ws.Range("AA341").Value = -1
$ws.Range("AB341").Value = -1
$ws.Range("AC341").Value = 0.95

# --- Row 342 becomes the former row-343 match (Club Celaya vs Oaxaca), ---
# --- with a couple of its odds adjusted, and old row 343 is removed.   ---
$ws.Range("B342").Value = 7641680
$ws.Range("E342").Value = 45347.83680555555
$ws.Range("F342").Value = "Club Celaya"
$ws.Range("G342").Value = "Oaxaca"
$ws.Range("K342").Value = 1.333
$ws.Range("L342").Value = 4.75
$ws.Range("M342").Value = 7.5
$ws.Range("N342").Value = 1.363
$ws.Range("O342").Value = 5
$ws.Range("P342").Value = 8
$ws.Range("Q342").Value = -1.5
$ws.Range("R342").Value = 1.95
$ws.Range("S342").Value = 1.85
$ws.Range("U342").Value = 1.95
$ws.Range("V342").Value = 1.85

# Drop the now-duplicated row 343 (its data was folded into row 342 above).
$ws.Range("A343").EntireRow.Delete()
